$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended to the "NEW" sheet (rows 76 and 77).
# Columns: A Caso, B F.De Reclamo, C Direccion, D Comuna, E OT, F Proveedor Asignado,
#          G Estado, H Observaciones, I Attachments, J Tipo de tarea, K Equipo,
#          L Tipo de Elemento, M Coordenada_X, N Coordenada_Y, O Operacion, P Zona, Q PD, R N2

$rows = @(
    @{
        RowNum = 76
        A = "7648"; B = "10/28/2025"; C = "JURAMENTO 5211"; D = "12"; E = "810458894"
        F = "NEW"; G = "Pendiente"; H = "Picada"; I = 1; J = "Cambio"; K = "Sin equipos"
        L = "Pasante"; M = -58.484108; N = -34.579014; O = "Paternal"; P = "Capital Norte"
        Q = "ATH-D"; R = "Fuera de Poligono OVL"
    },
    @{
        RowNum = 77
        A = "7663"; B = "10/28/2025"; C = "LAFINUR 2904"; D = "14"; E = "810458888"
        F = "NEW"; G = "Pendiente"; H = "Cambiar y reparar rienda"; I = 1; J = "Cambio"; K = "Sin equipos"
        L = "Terminal"; M = -58.415184; N = -34.581516; O = "Palermo"; P = "Capital Sur"
        Q = "AGU-M"; R = "Fuera de Poligono OVL"
    }
)

# Text-typed columns (must stay text even though some values look numeric/date)
$textCols = @("A", "B", "C", "D", "E", "F", "G", "H", "J", "K", "L", "O", "P", "Q", "R")
# Numeric-typed columns
$numCols = @("I", "M", "N")

foreach ($row in $rows) {
    $r = $row.RowNum

    # Force text number-format on the text columns of this row BEFORE assigning
    # values, so Excel doesn't auto-coerce numeric/date-looking strings
    # ("7648", "10/28/2025", "12", ...) into real numbers/dates.
    $textRange = $ws.Range("A$r" + ":" + "H$r")
    $textRange.NumberFormat = "@"
    $textRange2 = $ws.Range("J$r" + ":" + "L$r")
    $textRange2.NumberFormat = "@"
    $textRange3 = $ws.Range("O$r" + ":" + "R$r")
    $textRange3.NumberFormat = "@"

    foreach ($col in $textCols) {
        $ws.Range("$col$r").Value = $row[$col]
    }
    foreach ($col in $numCols) {
        $ws.Range("$col$r").Value = $row[$col]
    }

    # Reset the styling back to the default "Normal" style (no explicit style
    # index on these cells), matching the rest of the sheet's plain data rows.
    $ws.Range("A$r" + ":" + "R$r").Style = "Normal"
}

Write-Output "Added rows 76-77"
